# Automatische test-sync: 2025-06-29 14:13:50
# Appends a new "Testmail #2" log entry (row 11) to the Logs sheet and bumps
# the matching "Bestelling / Levering" tally on the Dashboard sheet.

$wb   = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

$newRow = 11

$logs.Range("A$newRow").Value = "Kun je 10 dozen schroeven bestellen?"
$logs.Range("B$newRow").Value = "mailmind.test@zohomail.eu"
$logs.Range("C$newRow").Value = "Testmail #2: Kun je 10 dozen schroeven bestellen?"
$logs.Range("D$newRow").Value = "Bestelling / Levering"
$logs.Range("E$newRow").Value = "Geachte klant," + [char]10 + `
    "Bedankt voor uw e-mail. Helaas kan ik geen bestellingen plaatsen, maar ik kan u doorverwijzen naar het bestelteam binnen ons bedrijf. Graag ontvang ik de contactgegevens van uw bedrijf, zodat ik de juiste persoon met u in contact kan brengen." + [char]10 + `
    "Ik zie uw reactie graag tegemoet." + [char]10 + `
    "Met vriendelijke groet," + [char]10 + `
    "[Naam]" + [char]10 + `
    "E-mailassistent"
$logs.Range("F$newRow").Value = "2025-06-29 14:13:30"
$logs.Range("G$newRow").Value = "Ja"
$logs.Range("H$newRow").Value = "Ja"
$logs.Range("I$newRow").Value = "Nee"

# Writing a multi-line value makes the host re-measure the row and pin an
# explicit height; AutoFit clears that pinned/custom height again so the new
# row matches the un-pinned rows above it.
$logs.Range("A$newRow" + ":I$newRow").EntireRow.AutoFit()

# Grow each conditional-formatting block so it also covers the new row.
# Modifying any one rule's applied range re-targets every cfRule sharing that
# block (they were all defined together over the same sqref).
$logs.Range("D2:D10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D$newRow"))
$logs.Range("G2:G10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G$newRow"))
$logs.Range("H2:H10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H$newRow"))
$logs.Range("I2:I10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I$newRow"))

# Dashboard tally for "Bestelling / Levering" goes from 2 to 3.
$dash.Range("B3").Value = 3
